# VACCINECERT-1633 Fixed CSV upload templates
#
# The sample row's "sampleDate" value (cell E2) held an incorrect serial
# date. Update it to the corrected date, then leave the active selection
# where the author's saved file shows it (G3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E2: sampleDate example value, 2021-11-01 (44501) -> 2021-11-16 (44516)
$ws.Range("E2").Value = 44516

# Restore the cursor/selection position recorded in the saved workbook.
$ws.Range("G3").Select()
